$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "test"

$ws.Range("C2").Value = [double]"0.3514970892147371"
$ws.Range("D2").Value = [double]"20.78828542189635"
$ws.Range("E2").Value = [double]"30.87602603714937"
$ws.Range("F2").Value = [double]"10.79598591506581"
$ws.Range("G2").Value = [double]"4.451076184841955"
$ws.Range("H2").Value = [double]"0.627849872862693"
$ws.Range("I2").Value = [double]"32.03452314979756"
$ws.Range("J2").Value = [double]"0.06143202706538113"
$ws.Range("K2").Value = [double]"0.01325335131680187"
$ws.Range("M2").Value = [double]"2.952091034383236e-05"
$ws.Range("N2").Value = [double]"7.703922829993426e-07"
$ws.Range("O2").Value = [double]"1.985961376398509e-09"
$ws.Range("P2").Value = [double]"6.888519895734363e-15"
$ws.Range("Q2").Value = [double]"2.108854848703124e-15"
$ws.Range("R2").Value = [double]"1.589305646544744e-17"
$ws.Range("S2").Value = [double]"8.903787417131954e-19"
$ws.Range("T2").Value = [double]"4.527155162802647e-19"
$ws.Range("U2").Value = [double]"2.425948187943382e-22"
$ws.Range("V2").Value = [double]"4.502947776946262e-20"
$ws.Range("X2").Value = [double]"1.51513121912426e-22"
$ws.Range("Y2").Value = [double]"1.182070249297177e-11"
$ws.Range("Z2").Value = [double]"6.100914311238397e-32"
$ws.Range("AB2").Value = [double]"1.793662049747534e-09"
$ws.Range("AD2").Value = [double]"1.238707869787348e-09"
$ws.Range("AE2").Value = [double]"7.493738925596594e-28"
$ws.Range("AF2").Value = [double]"3.247262401926823e-11"
$ws.Range("AG2").Value = [double]"1.136394757033844e-12"
$ws.Range("AH2").Value = [double]"3.423032517393449e-05"
$ws.Range("AI2").Value = [double]"6.424093505802726e-06"
$ws.Range("AJ2").Value = [double]"4.223548609119147e-12"
$ws.Range("AL2").Value = [double]"4.402624789659949e-14"
$ws.Range("AM2").Value = [double]"1.625692340411667e-15"
$ws.Range("AN2").Value = [double]"1.658415652298465e-18"
$ws.Range("AR2").Value = [double]"0.0004986683431200219"
$ws.Range("AS2").Value = [double]"0.2396853194917016"
$ws.Range("AT2").Value = [double]"1.350941178606322"
$ws.Range("AU2").Value = [double]"1.849507820813773"
$ws.Range("AV2").Value = [double]"2.242118071798063"
$ws.Range("AW2").Value = [double]"0.9755572931588912"
$ws.Range("AX2").Value = [double]"93.1619981305147"
$ws.Range("AY2").Value = [double]"0.1315526878297018"
$ws.Range("AZ2").Value = [double]"0.04788721394350171"
$ws.Range("BB2").Value = [double]"0.0002403072808293013"
$ws.Range("BC2").Value = [double]"1.037861407585402e-05"
$ws.Range("BD2").Value = [double]"2.557471497725458e-08"
$ws.Range("BE2").Value = [double]"1.318618384030238e-13"
$ws.Range("BF2").Value = [double]"4.170766350116881e-14"
$ws.Range("BG2").Value = [double]"5.396620021262361e-16"
$ws.Range("BH2").Value = [double]"3.703919845779697e-17"
$ws.Range("BI2").Value = [double]"1.835911599244912e-17"
$ws.Range("BJ2").Value = [double]"1.10605388416395e-20"
$ws.Range("BK2").Value = [double]"1.868722920433879e-18"
$ws.Range("BL2").Value = [double]"1.372513943895427e-36"
$ws.Range("BM2").Value = [double]"6.66979854539755e-21"
$ws.Range("BN2").Value = [double]"2.352174912078328e-12"
$ws.Range("BO2").Value = [double]"1.226959087825827e-30"
$ws.Range("BP2").Value = [double]"2.525098215537728e-34"
$ws.Range("BQ2").Value = [double]"1.906161870997596e-10"
$ws.Range("BR2").Value = [double]"7.237344690264337e-36"
$ws.Range("BS2").Value = [double]"8.267894570785628e-11"
$ws.Range("BT2").Value = [double]"5.050392837545223e-26"
$ws.Range("BU2").Value = [double]"2.099413986548547e-12"
$ws.Range("BV2").Value = [double]"1.105959375060462e-13"
$ws.Range("BW2").Value = [double]"2.245291874126622e-06"
$ws.Range("BX2").Value = [double]"6.584598654044067e-07"
$ws.Range("BY2").Value = [double]"8.412388872631388e-13"
$ws.Range("CA2").Value = [double]"3.019798959613349e-15"
$ws.Range("CB2").Value = [double]"2.195392778918456e-16"
$ws.Range("CC2").Value = [double]"2.798791000157198e-19"
$ws.Range("CG2").Value = [double]"3.067418861684802e-07"
$ws.Range("CH2").Value = [double]"7.697561138950245e-05"
$ws.Range("CI2").Value = [double]"0.0003951969342515768"
$ws.Range("CJ2").Value = [double]"0.0007353375359049172"
$ws.Range("CK2").Value = [double]"0.001769693984767946"
$ws.Range("CL2").Value = [double]"0.01080905770883727"
$ws.Range("CM2").Value = [double]"99.90194156733509"
$ws.Range("CN2").Value = [double]"0.02303541164484406"
$ws.Range("CO2").Value = [double]"0.05846258036168871"
$ws.Range("CQ2").Value = [double]"0.001974843321152549"
$ws.Range("CR2").Value = [double]"0.0007979963431563413"
$ws.Range("CS2").Value = [double]"1.032401008978719e-06"
$ws.Range("CT2").Value = [double]"5.460766637238191e-11"
$ws.Range("CU2").Value = [double]"1.912915197662088e-11"
$ws.Range("CV2").Value = [double]"1.92853928101543e-12"
$ws.Range("CW2").Value = [double]"1.918191780804434e-13"
$ws.Range("CX2").Value = [double]"1.26890951443159e-13"
$ws.Range("CY2").Value = [double]"1.875318740125801e-16"
$ws.Range("CZ2").Value = [double]"7.960450102945499e-15"
$ws.Range("DA2").Value = [double]"1.75134270280955e-24"
$ws.Range("DB2").Value = [double]"9.096858182708633e-17"
$ws.Range("DC2").Value = [double]"1.382801702074209e-27"
$ws.Range("DD2").Value = [double]"7.009272411816171e-24"
$ws.Range("DE2").Value = [double]"1.387040749577175e-24"
$ws.Range("DF2").Value = [double]"4.256969890419802e-28"
$ws.Range("DG2").Value = [double]"5.582326328030815e-25"
$ws.Range("DH2").Value = [double]"3.014638930135636e-35"
$ws.Range("DI2").Value = [double]"4.227184969208193e-21"
$ws.Range("DJ2").Value = [double]"3.820873131077883e-14"
$ws.Range("DK2").Value = [double]"1.804401134676397e-15"
$ws.Range("DM2").Value = [double]"1.724735700785842e-31"
$ws.Range("DN2").Value = [double]"6.558189004058982e-15"
$ws.Range("DP2").Value = [double]"3.917137643429477e-17"
$ws.Range("DQ2").Value = [double]"1.73714162077538e-18"
$ws.Range("DR2").Value = [double]"1.295356152450408e-21"
$ws.Range("DV2").Value = [double]"1.528277188631658e-09"
$ws.Range("DW2").Value = [double]"6.100385015861622e-08"
$ws.Range("DX2").Value = [double]"1.335318905772428e-07"
$ws.Range("DY2").Value = [double]"1.149051751719298e-27"
$ws.Range("DZ2").Value = [double]"5.390242555161386e-20"
$ws.Range("EA2").Value = [double]"1.173198293448328e-12"
$ws.Range("EB2").Value = [double]"0.004361323251807471"
$ws.Range("EC2").Value = [double]"4.091998100005618e-07"
$ws.Range("ED2").Value = [double]"0.0001601297257438549"
$ws.Range("EF2").Value = [double]"0.02107939209528966"
$ws.Range("EG2").Value = [double]"95.81310444741889"
$ws.Range("EH2").Value = [double]"0.01205497480889129"
$ws.Range("EI2").Value = [double]"0.008981639649232617"
$ws.Range("EJ2").Value = [double]"0.007384948091595608"
$ws.Range("EK2").Value = [double]"1.844469687512324"
$ws.Range("EL2").Value = [double]"0.8465975968281253"
$ws.Range("EM2").Value = [double]"1.302965180821242"
$ws.Range("EN2").Value = [double]"0.07761561717450505"
$ws.Range("EO2").Value = [double]"0.01614678819581189"
$ws.Range("EP2").Value = [double]"3.029070563584932e-08"
$ws.Range("EQ2").Value = [double]"0.03574249546428865"
$ws.Range("ER2").Value = [double]"5.047492972074742e-09"
$ws.Range("ES2").Value = [double]"1.461867202229782e-06"
$ws.Range("ET2").Value = [double]"5.460348203451919e-07"
$ws.Range("EU2").Value = [double]"4.386567911546125e-07"
$ws.Range("EV2").Value = [double]"2.361146573277556e-07"
$ws.Range("EW2").Value = [double]"2.610963785359511e-07"
$ws.Range("EX2").Value = [double]"0.002892995292462147"
$ws.Range("EY2").Value = [double]"5.752873497525722e-09"
$ws.Range("EZ2").Value = [double]"2.461432399691903e-10"
$ws.Range("FA2").Value = [double]"0.005132688645593005"
$ws.Range("FB2").Value = [double]"0.00130650369525937"
$ws.Range("FC2").Value = [double]"9.465313736330314e-10"
$ws.Range("FE2").Value = [double]"9.888029907146588e-12"
$ws.Range("FF2").Value = [double]"4.760742266249713e-13"
$ws.Range("FG2").Value = [double]"5.119282742198963e-16"
